$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add final "Post Treatment" outcome measure values for the parent data (column D)
$values = @(45, 49, 46, 65, 48, 56, 61, 59, 53, 47, 46, 54, 53)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Update the active selection to reflect where entry left off
$ws.Range("D15").Select() | Out-Null
